$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.142.11'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.596.30'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9993'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '302.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3779'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3594'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.97'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.255'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9992'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08117'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.565'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001247'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.353'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("D17").Value = '1.593.96'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06835'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.512'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("D24").Value = '23.148.23'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.394'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.947'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.228'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.363'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.768'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("D33").Value = '1.766.63'
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9712'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07533'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02694'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.138'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08796'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2492'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.15%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7106'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.363'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6497'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.008'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.284'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07910'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.206'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.215'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.47%  '
